$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 149, 150
$ws.Range("B149").Value = 63902
$ws.Range("D149").Value = 32.02
$ws.Range("E149").Value = 34.04
$ws.Range("F149").Value = 2
$ws.Range("G149").Value = 64.04
$ws.Range("B150").Value = 48654
$ws.Range("D150").Value = 32.02
$ws.Range("E150").Value = 38.26
$ws.Range("F150").Value = -1
$ws.Range("G150").Value = -32.02

# Rows 183, 184
$ws.Range("B183").Value = 64329
$ws.Range("D183").Value = 120.69
$ws.Range("E183").Value = 128.32
$ws.Range("F183").Value = 6
$ws.Range("G183").Value = 724.14
$ws.Range("B184").Value = 57552
$ws.Range("D184").Value = 120.69
$ws.Range("E184").Value = 136.86
$ws.Range("F184").Value = -5
$ws.Range("G184").Value = -603.45

# Rows 279, 280
$ws.Range("B279").Value = 64973
$ws.Range("D279").Value = 33.3
$ws.Range("E279").Value = 35.4
$ws.Range("F279").Value = 150
$ws.Range("G279").Value = 4995
$ws.Range("B280").Value = 48706
$ws.Range("D280").Value = 33.3
$ws.Range("E280").Value = 39.8
$ws.Range("F280").Value = -144
$ws.Range("G280").Value = -4795.2

# Rows 313, 314
$ws.Range("B313").Value = 57854
$ws.Range("D313").Value = 305.84
$ws.Range("E313").Value = 325.16
$ws.Range("F313").Value = 2
$ws.Range("G313").Value = 611.68
$ws.Range("B314").Value = 62997
$ws.Range("D314").Value = 305.84
$ws.Range("E314").Value = 325.16
$ws.Range("F314").Value = 72
$ws.Range("G314").Value = 22020.48

# Rows 316, 317, 318
$ws.Range("B316").Value = 63565
$ws.Range("D316").Value = 102.71
$ws.Range("E316").Value = 109.19
$ws.Range("F316").Value = 60
$ws.Range("G316").Value = 6162.6
$ws.Range("B317").Value = 57077
$ws.Range("D317").Value = 93.08
$ws.Range("E317").Value = 111.2
$ws.Range("F317").Value = 1
$ws.Range("G317").Value = 93.08
$ws.Range("B318").Value = 61610
$ws.Range("D318").Value = 102.71
$ws.Range("E318").Value = 122.71
$ws.Range("F318").Value = -58
$ws.Range("G318").Value = -5957.18

# Rows 350, 351, 352
$ws.Range("B350").Value = 63571
$ws.Range("D350").Value = 143.48
$ws.Range("E350").Value = 152.53
$ws.Range("F350").Value = 27
$ws.Range("G350").Value = 3873.96
$ws.Range("B351").Value = 57802
$ws.Range("D351").Value = 143.48
$ws.Range("E351").Value = 162.71
$ws.Range("F351").Value = -79
$ws.Range("G351").Value = -11334.92
$ws.Range("B352").Value = 63531
$ws.Range("D352").Value = 143.48
$ws.Range("E352").Value = 152.53
$ws.Range("F352").Value = 80
$ws.Range("G352").Value = 11478.4

# Rows 379, 380
$ws.Range("B379").Value = 61608
$ws.Range("D379").Value = 129.01
$ws.Range("E379").Value = 154.12
$ws.Range("F379").Value = -56
$ws.Range("G379").Value = -7224.56
$ws.Range("B380").Value = 63564
$ws.Range("D380").Value = 129.01
$ws.Range("E380").Value = 137.16
$ws.Range("F380").Value = 57
$ws.Range("G380").Value = 7353.57

# Rows 382, 383
$ws.Range("B382").Value = 63560
$ws.Range("D382").Value = 126.86
$ws.Range("E382").Value = 134.87
$ws.Range("F382").Value = 104
$ws.Range("G382").Value = 13193.44
$ws.Range("B383").Value = 60325
$ws.Range("D383").Value = 126.86
$ws.Range("E383").Value = 151.57
$ws.Range("F383").Value = -102
$ws.Range("G383").Value = -12939.72

# Rows 400, 401
$ws.Range("B400").Value = 57835
$ws.Range("D400").Value = 59.13
$ws.Range("E400").Value = 70.65
$ws.Range("F400").Value = 1
$ws.Range("G400").Value = 59.13
$ws.Range("B401").Value = 62933
$ws.Range("D401").Value = 59.13
$ws.Range("E401").Value = 70.65
$ws.Range("F401").Value = 146
$ws.Range("G401").Value = 8632.98

# Rows 421, 422
$ws.Range("B421").Value = 57857
$ws.Range("D421").Value = 151.17
$ws.Range("E421").Value = 180.62
$ws.Range("F421").Value = 3
$ws.Range("G421").Value = 453.51
$ws.Range("B422").Value = 63008
$ws.Range("D422").Value = 151.17
$ws.Range("E422").Value = 180.62
$ws.Range("F422").Value = 504
$ws.Range("G422").Value = 76189.68

# Rows 579, 580
$ws.Range("B579").Value = 53757
$ws.Range("D579").Value = 13.45
$ws.Range("E579").Value = 16.08
$ws.Range("F579").Value = -159
$ws.Range("G579").Value = -2138.55
$ws.Range("B580").Value = 65069
$ws.Range("D580").Value = 13.45
$ws.Range("E580").Value = 14.3
$ws.Range("F580").Value = 172
$ws.Range("G580").Value = 2313.4

# Rows 581, 582
$ws.Range("B581").Value = 53602
$ws.Range("D581").Value = 13.15
$ws.Range("E581").Value = 15.69
$ws.Range("F581").Value = -231
$ws.Range("G581").Value = -3037.65
$ws.Range("B582").Value = 65068
$ws.Range("D582").Value = 13.15
$ws.Range("E582").Value = 13.97
$ws.Range("F582").Value = 232
$ws.Range("G582").Value = 3050.8

# Rows 583, 584
$ws.Range("B583").Value = 65066
$ws.Range("D583").Value = 12.81
$ws.Range("E583").Value = 13.61
$ws.Range("F583").Value = 313
$ws.Range("G583").Value = 4009.53
$ws.Range("B584").Value = 53263
$ws.Range("D584").Value = 12.81
$ws.Range("E584").Value = 15.29
$ws.Range("F584").Value = -309
$ws.Range("G584").Value = -3958.29

# Rows 586, 587
$ws.Range("B586").Value = 45695
$ws.Range("D586").Value = 19.73
$ws.Range("E586").Value = 23.58
$ws.Range("F586").Value = -36
$ws.Range("G586").Value = -710.28
$ws.Range("B587").Value = 64915
$ws.Range("D587").Value = 19.73
$ws.Range("E587").Value = 20.98
$ws.Range("F587").Value = 40
$ws.Range("G587").Value = 789.2

# Rows 593, 594
$ws.Range("B593").Value = 64927
$ws.Range("D593").Value = 16.22
$ws.Range("E593").Value = 17.26
$ws.Range("F593").Value = 295
$ws.Range("G593").Value = 4784.9
$ws.Range("B594").Value = 45718
$ws.Range("D594").Value = 16.22
$ws.Range("E594").Value = 19.38
$ws.Range("F594").Value = -294
$ws.Range("G594").Value = -4768.68

# Rows 601, 602
$ws.Range("B601").Value = 45702
$ws.Range("D601").Value = 26.3
$ws.Range("E601").Value = 31.43
$ws.Range("F601").Value = -215
$ws.Range("G601").Value = -5654.5
$ws.Range("B602").Value = 64919
$ws.Range("D602").Value = 26.3
$ws.Range("E602").Value = 27.97
$ws.Range("F602").Value = 224
$ws.Range("G602").Value = 5891.2
